# EI Variable Installments T1 scenarios
# Applies the recorded value/selection/column-width edits to the three
# affected worksheets: Summary, Repayment schedule, Transactions.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Cells.Item(2, 2).Value = 836.76    # B2
$wsSummary.Cells.Item(2, 5).Value = 9163.24   # E2
$wsSummary.Cells.Item(2, 6).Value = 849.4     # F2

$wsSummary.Cells.Item(3, 1).Value = 561.21    # A3
$wsSummary.Cells.Item(3, 5).Value = 510.25    # E3
$wsSummary.Cells.Item(3, 6).Value = 38.32     # F3

$wsSummary.Cells.Item(5, 1).Value = 200       # A5
$wsSummary.Cells.Item(5, 2).Value = 100       # B5
$wsSummary.Cells.Item(5, 5).Value = 100       # E5
$wsSummary.Cells.Item(5, 6).Value = 100       # F5

[void]$wsSummary.Activate()
[void]$wsSummary.Range("A7:XFD15").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsSched = $wb.Worksheets.Item("Repayment schedule")

$wsSched.Cells.Item(3, 10).Value = 100        # J3

$wsSched.Cells.Item(5, 6).Value = 849.4       # F5
$wsSched.Cells.Item(5, 7).Value = 8313.84     # G5
$wsSched.Cells.Item(5, 8).Value = 38.32       # H5
$wsSched.Cells.Item(5, 10).Value = 100        # J5
$wsSched.Cells.Item(5, 11).Value = 987.72     # K5
$wsSched.Cells.Item(5, 16).Value = 987.72     # P5

$wsSched.Cells.Item(6, 6).Value = 794.33      # F6
$wsSched.Cells.Item(6, 7).Value = 7519.51     # G6
$wsSched.Cells.Item(6, 8).Value = 93.39       # H6

$wsSched.Cells.Item(7, 6).Value = 813.55      # F7
$wsSched.Cells.Item(7, 7).Value = 6705.96     # G7
$wsSched.Cells.Item(7, 8).Value = 74.17       # H7

$wsSched.Cells.Item(8, 6).Value = 819.37      # F8
$wsSched.Cells.Item(8, 7).Value = 5886.59     # G8
$wsSched.Cells.Item(8, 8).Value = 68.35       # H8

$wsSched.Cells.Item(9, 6).Value = 829.66      # F9
$wsSched.Cells.Item(9, 7).Value = 5056.93     # G9
$wsSched.Cells.Item(9, 8).Value = 58.06       # H9

$wsSched.Cells.Item(10, 6).Value = 836.18     # F10
$wsSched.Cells.Item(10, 7).Value = 4220.75    # G10
$wsSched.Cells.Item(10, 8).Value = 51.54      # H10

$wsSched.Cells.Item(11, 6).Value = 844.7      # F11
$wsSched.Cells.Item(11, 7).Value = 3376.05    # G11
$wsSched.Cells.Item(11, 8).Value = 43.02      # H11

$wsSched.Cells.Item(12, 6).Value = 854.42     # F12
$wsSched.Cells.Item(12, 7).Value = 2521.63    # G12
$wsSched.Cells.Item(12, 8).Value = 33.3       # H12

$wsSched.Cells.Item(13, 6).Value = 862.02     # F13
$wsSched.Cells.Item(13, 7).Value = 1659.61    # G13
$wsSched.Cells.Item(13, 8).Value = 25.7       # H13

$wsSched.Cells.Item(14, 6).Value = 871.35     # F14
$wsSched.Cells.Item(14, 7).Value = 788.26     # G14
$wsSched.Cells.Item(14, 8).Value = 16.37      # H14

$wsSched.Cells.Item(15, 6).Value = 788.26     # F15
$wsSched.Cells.Item(15, 8).Value = 8.03       # H15
$wsSched.Cells.Item(15, 11).Value = 796.29    # K15
$wsSched.Cells.Item(15, 16).Value = 796.29    # P15

[void]$wsSched.Activate()
[void]$wsSched.Range("F20").Select()

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$wsTx = $wb.Worksheets.Item("Transactions")

$wsTx.Cells.Item(2, 1).Value = 191            # A2
$wsTx.Cells.Item(2, 10).Value = 9163.24       # J2
$wsTx.Cells.Item(2, 10).NumberFormat = "#,##0.00"

$wsTx.Cells.Item(3, 1).Value = 189            # A3
$wsTx.Cells.Item(3, 6).Value = 836.76         # F3
$wsTx.Cells.Item(3, 9).Value = 100            # I3
$wsTx.Cells.Item(3, 10).Value = 4163.24       # J3
$wsTx.Cells.Item(3, 10).NumberFormat = "#,##0.00"

$wsTx.Cells.Item(4, 1).Value = 180            # A4

$wsTx.Columns.Item(1).ColumnWidth = 3.1666666666666665

[void]$wsTx.Activate()
[void]$wsTx.Range("C13").Select()
